# Applies the "Automatic update of files." commit:
#  - Column C ("Förändrad") advances by one day (45178 -> 45179) for every
#    data row (rows 2-319) on the "Avverkningsanmälningar" sheet.
#  - Row 3 additionally loses one "NT" hit, one "Rödlistade" hit, one "Alla
#    arter" hit, and the "Gul taggsvamp" line in its species list (R3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldDate = 45178
$newDate = 45179

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 319
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq $oldDate) {
        $cell.Value2 = $newDate
    }
}

# Row-3 specific count adjustments.
$ws.Range("J3").Value2 = 1
$ws.Range("O3").Value2 = 2
$ws.Range("Q3").Value2 = 4

# Row-3 species list: drop the "Gul taggsvamp" line.
$nl = "`r`n"
$ws.Range("R3").Value2 = "Knärot" + $nl + "Skirmossa" + $nl + "Havstulpanlav" + $nl + "Korallblylav"
